# Weekly refresh of Fruta/Hortaliza (Chirimoya) price rows: the source
# rows for this market were re-pulled and land on different rows than
# before. Re-write the per-row fields (Fecha, Calidad, Volumen, prices)
# so every row ends up holding the data for its new source row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => Fecha (serial date), Calidad, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Precio $/Kg
$rows = @(
    @{ Row = 2;  D = 44841; L = "Primera";  M = 60;  N = 23000; O = 24000; P = 23500; S = 2350 },
    @{ Row = 3;  D = 44848; L = "Especial"; M = 60;  N = 24000; O = 25000; P = 24500; S = 2450 },
    @{ Row = 4;  D = 44848; L = "Primera";  M = 120; N = 21000; O = 22000; P = 21500; S = 2150 },
    @{ Row = 5;  D = 44447; L = "Primera";  M = 60;  N = 21000; O = 22000; P = 21500; S = 2150 },
    @{ Row = 6;  D = 44461; L = "Especial"; M = 60;  N = 31000; O = 32000; P = 31500; S = 3150 },
    @{ Row = 7;  D = 44461; L = "Primera";  M = 30;  N = 30000; O = 30000; P = 30000; S = 3000 },
    @{ Row = 8;  D = 44460; L = "Especial"; M = 60;  N = 31000; O = 32000; P = 31500; S = 3150 },
    @{ Row = 9;  D = 44460; L = "Primera";  M = 30;  N = 30000; O = 30000; P = 30000; S = 3000 },
    @{ Row = 10; D = 44874; L = "Especial"; M = 30;  N = 25000; O = 25000; P = 25000; S = 2500 },
    @{ Row = 11; D = 44874; L = "Primera";  M = 80;  N = 23000; O = 24000; P = 23500; S = 2350 },
    @{ Row = 13; D = 44446; L = "Primera";  M = 60;  N = 21000; O = 22000; P = 21500; S = 2150 },
    @{ Row = 14; D = 44839; L = "Primera";  M = 120; N = 25000; O = 26000; P = 25500; S = 2550 },
    @{ Row = 16; D = 44868; L = "Especial"; M = 60;  N = 26000; O = 26000; P = 26000; S = 2600 },
    @{ Row = 17; D = 44487; L = "Primera";  M = 30;  N = 23000; O = 24000; P = 23500; S = 2350 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value = $r.D    # D: Fecha
    $ws.Cells.Item($row, 12).Value = $r.L   # L: Calidad
    $ws.Cells.Item($row, 13).Value = $r.M   # M: Volumen
    $ws.Cells.Item($row, 14).Value = $r.N   # N: Precio minimo
    $ws.Cells.Item($row, 15).Value = $r.O   # O: Precio maximo
    $ws.Cells.Item($row, 16).Value = $r.P   # P: Precio promedio ponderado
    $ws.Cells.Item($row, 19).Value = $r.S   # S: Precio $/Kg
}
